# Insert a new weekly data row at row 43 (pushing existing rows 43-55 down to 44-56),
# matching the "Fruta / hortaliza, semanal" update that adds a new price observation
# dated 2021-11-24 (serial 44524) for Macroferia Regional de Talca - Espárragos.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 43; this shifts formatting/values of the
# old rows 43..55 down to 44..56 automatically, same as Excel's UI "Insert Row".
$ws.Rows.Item(43).Insert()

# Populate the newly inserted row 43 with the new record's values.
$ws.Cells.Item(43, 1).Value2 = 5
$ws.Cells.Item(43, 2).Value2 = "Macroferia Regional de Talca"
$ws.Cells.Item(43, 3).Value2 = "Maule"
$ws.Cells.Item(43, 4).Value2 = 44524
$ws.Cells.Item(43, 5).Value2 = 7
$ws.Cells.Item(43, 6).Value2 = 300000000
$ws.Cells.Item(43, 7).Value2 = "Espárragos"
$ws.Cells.Item(43, 8).Value2 = "Verde"
$ws.Cells.Item(43, 9).Value2 = "Primera"
$ws.Cells.Item(43, 10).Value2 = 3000
$ws.Cells.Item(43, 11).Value2 = 1000
$ws.Cells.Item(43, 12).Value2 = 1000
$ws.Cells.Item(43, 13).Value2 = 1000
$ws.Cells.Item(43, 14).Value2 = "`$/kilo"
$ws.Cells.Item(43, 15).Value2 = "Región del Maule"
$ws.Cells.Item(43, 16).Value2 = 1000
$ws.Cells.Item(43, 17).Value2 = 1
$ws.Cells.Item(43, 18).Value2 = "Hortaliza"

# Make sure the date cell keeps the same number format as the other date cells
# in column D (style index 2 in styles.xml => numFmtId 165, "YYYY-MM-DD HH:MM:SS").
$ws.Cells.Item(43, 4).NumberFormat = $ws.Cells.Item(44, 4).NumberFormat
